$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting the existing rows 57-59 down to 58-60.
$ws.Rows.Item(57).Insert()

# Populate the newly-inserted row 57 with its values.
$ws.Cells.Item(57, 1).Value = 7
$ws.Cells.Item(57, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(57, 3).Value = "Ñuble"
$ws.Cells.Item(57, 4).Value = 45212
$ws.Cells.Item(57, 4).Style = $ws.Cells.Item(58, 4).Style
$ws.Cells.Item(57, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(57, 6).Value = 300000000
$ws.Cells.Item(57, 7).Value = "Espárragos"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 600
$ws.Cells.Item(57, 11).Value = 1200
$ws.Cells.Item(57, 12).Value = 1300
$ws.Cells.Item(57, 13).Value = 1250
$ws.Cells.Item(57, 14).Value = "$/kilo"
$ws.Cells.Item(57, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(57, 16).Value = 1250
$ws.Cells.Item(57, 17).Value = 1
$ws.Cells.Item(57, 18).Value = "Hortaliza"
